# Update "想去人数" (want-to-go count) values in column F
# for sheet "展览" (rows 2-24) and sheet "全部类型" (rows 2-26).
# Both sheets list the same con/event entries (全部类型 has a couple of
# extra rows from the 演出 sheet interleaved), so the row numbers differ
# slightly between the two sheets but the target values are identical.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new value, for sheet "展览"
$exhibitUpdates = @{
    2  = 477
    4  = 8022
    13 = 457
    14 = 70
    17 = 5916
    18 = 189
    20 = 1927
    21 = 20
    22 = 32
    24 = 410
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row -> new value, for sheet "全部类型"
$allTypesUpdates = @{
    2  = 477
    4  = 8022
    13 = 457
    14 = 70
    18 = 5916
    20 = 189
    22 = 1927
    23 = 20
    24 = 32
    26 = 410
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allTypesUpdates[$row]
}

$wb.Save()
